$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Test Steps": insert 5 new rows (My Account page verification steps)
# right before the "Menu_Navigation" block, and renumber the TS ids of the
# "Logout_01" rows that follow.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Test Steps")

# Insert 5 blank rows before current row 13 (Menu_Navigation block)
$ws1.Range("A13:G17").Insert()

# Copy the formatting (borders/styles) of an existing data row into the
# newly inserted rows.
$ws1.Range("A12:G12").Copy()
$ws1.Range("A13:G17").PasteSpecial(-4122)

# Row 13
$ws1.Range("A13").Value = "Verify_MyAccount_Page"
$ws1.Range("B13").Value = "TS_014"
$ws1.Range("C13").Value = "Verify the element btn_OrderHistory of My Account Page"
$ws1.Range("D13").Value = "btn_OrderHistory"
$ws1.Range("E13").Value = "verifyElement"
$ws1.Range("G13").Value = "PASS"

# Row 14
$ws1.Range("A14").Value = "Verify_MyAccount_Page"
$ws1.Range("B14").Value = "TS_015"
$ws1.Range("C14").Value = "Verify the element btn_MycreditSlips of My Account Page"
$ws1.Range("D14").Value = "btn_MycreditSlips"
$ws1.Range("E14").Value = "verifyElement"
$ws1.Range("G14").Value = "PASS"

# Row 15
$ws1.Range("A15").Value = "Verify_MyAccount_Page"
$ws1.Range("B15").Value = "TS_016"
$ws1.Range("C15").Value = "Verify the element btn_Myaddress of My Account Page"
$ws1.Range("D15").Value = "btn_Myaddress"
$ws1.Range("E15").Value = "verifyElement"
$ws1.Range("G15").Value = "PASS"

# Row 16
$ws1.Range("A16").Value = "Verify_MyAccount_Page"
$ws1.Range("B16").Value = "TS_017"
$ws1.Range("C16").Value = "Verify the element btn_Mypersonalinfo of My Account Page"
$ws1.Range("D16").Value = "btn_Mypersonalinfo"
$ws1.Range("E16").Value = "verifyElement"
$ws1.Range("G16").Value = "PASS"

# Row 17
$ws1.Range("A17").Value = "Verify_MyAccount_Page"
$ws1.Range("B17").Value = "TS_018"
$ws1.Range("C17").Value = "Verify the elements btn_Mywhishlist of My Account Page"
$ws1.Range("D17").Value = "btn_Mywhishlist"
$ws1.Range("E17").Value = "verifyElement"
$ws1.Range("G17").Value = "PASS"

# The former rows 13-16 (Menu_Navigation x2, Logout_01 x2) are now rows
# 18-21. The Logout_01 rows' TS ids are renumbered to TS_019 / TS_020.
$ws1.Range("B20").Value = "TS_019"
$ws1.Range("B21").Value = "TS_020"

# Sheet cosmetics: widen column A and move the active selection.
$ws1.Columns.Item(1).ColumnWidth = 20.333333333333336
$ws1.Range("F2").Select()

# ---------------------------------------------------------------------------
# Sheet "Test Cases": insert a new row for the "Verify_MyAccount_Page" case
# right before the "Menu_Navigation" row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Test Cases")

$ws2.Range("A4:D4").Insert()
$ws2.Range("A3:D3").Copy()
$ws2.Range("A4:D4").PasteSpecial(-4122)

$ws2.Range("A4").Value = "Verify_MyAccount_Page"
$ws2.Range("B4").Value = "Verify the elements of My Account Page"
$ws2.Range("C4").Value = "Yes"
$ws2.Range("D4").Value = "PASS"

$ws2.Columns.Item(1).ColumnWidth = 20.333333333333336
$ws2.Range("C7").Select()

$ws1.Activate()
